# Update NATMI LR-pair edge statistics following Dr Hou advice
# (recomputed with updated ligand/receptor-expressing cell counts: 1 -> 3)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 14.768619
$ws.Range("H2").Value = 44.305857
$ws.Range("I2").Value = 0.9736679609684162
$ws.Range("J2").Value = 0.9736679609684162
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7530926666666667
$ws.Range("N2").Value = 2.259278
$ws.Range("O2").Value = 0.1998288175720169
$ws.Range("P2").Value = 0.1998288175720169
$ws.Range("Q2").Value = 11.122138665694
$ws.Range("R2").Value = 100.099247991246
$ws.Range("S2").Value = 0.1945669173480753
$ws.Range("T2").Value = 0.1945669173480753
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 14.768619
$ws.Range("H3").Value = 44.305857
$ws.Range("I3").Value = 0.9736679609684162
$ws.Range("J3").Value = 0.9736679609684162
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.557543666666667
$ws.Range("N3").Value = 4.672631
$ws.Range("O3").Value = 0.4132852741806677
$ws.Range("P3").Value = 0.4132852741806678
$ws.Range("Q3").Value = 23.002768988863
$ws.Range("R3").Value = 207.024920899767
$ws.Range("S3").Value = 0.4024026302097636
$ws.Range("T3").Value = 0.4024026302097636
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 14.768619
$ws.Range("H4").Value = 44.305857
$ws.Range("I4").Value = 0.9736679609684162
$ws.Range("J4").Value = 0.9736679609684162
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.458052666666666
$ws.Range("N4").Value = 4.374158
$ws.Range("O4").Value = 0.3868859082473153
$ws.Range("P4").Value = 0.3868859082473153
$ws.Range("Q4").Value = 21.533424315934
$ws.Range("R4").Value = 193.800818843406
$ws.Range("S4").Value = 0.3766984134105772
$ws.Range("T4").Value = 0.3766984134105773
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.399405
$ws.Range("H5").Value = 1.198215
$ws.Range("I5").Value = 0.02633203903158381
$ws.Range("J5").Value = 0.02633203903158381
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7530926666666667
$ws.Range("N5").Value = 2.259278
$ws.Range("O5").Value = 0.1998288175720169
$ws.Range("P5").Value = 0.1998288175720169
$ws.Range("Q5").Value = 0.30078897653
$ws.Range("R5").Value = 2.70710078877
$ws.Range("S5").Value = 0.005261900223941589
$ws.Range("T5").Value = 0.00526190022394159
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.399405
$ws.Range("H6").Value = 1.198215
$ws.Range("I6").Value = 0.02633203903158381
$ws.Range("J6").Value = 0.02633203903158381
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.557543666666667
$ws.Range("N6").Value = 4.672631
$ws.Range("O6").Value = 0.4132852741806677
$ws.Range("P6").Value = 0.4132852741806678
$ws.Range("Q6").Value = 0.622090728185
$ws.Range("R6").Value = 5.598816553665
$ws.Range("S6").Value = 0.01088264397090416
$ws.Range("T6").Value = 0.01088264397090416
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.399405
$ws.Range("H7").Value = 1.198215
$ws.Range("I7").Value = 0.02633203903158381
$ws.Range("J7").Value = 0.02633203903158381
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.458052666666666
$ws.Range("N7").Value = 4.374158
$ws.Range("O7").Value = 0.3868859082473153
$ws.Range("P7").Value = 0.3868859082473153
$ws.Range("Q7").Value = 0.58235352533
$ws.Range("R7").Value = 5.24118172797
$ws.Range("S7").Value = 0.01018749483673806
$ws.Range("T7").Value = 0.01018749483673806